$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "SA"

# 2. Tiny correction of two existing values (last-digit rounding fix)
$ws.Cells.Item(13, 6).Value = 0.9938684977613513   # F13
$ws.Cells.Item(13, 16).Value = 0.9942635674110872  # P13

# 3. Append a new row (row 16) with data, mirroring the style/format of row 15
$ws.Range("A15:P15").Copy()
$ws.Range("A16:P16").PasteSpecial(-4122)

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(16, 3).Value = 1.052557562715659
$ws.Cells.Item(16, 4).Value = 0.9131717376341408
$ws.Cells.Item(16, 5).Value = 1.016609016762428
$ws.Cells.Item(16, 6).Value = 0.9761360620500866
$ws.Cells.Item(16, 7).Value = 1.052557562715659
$ws.Cells.Item(16, 8).Value = 0.9131717376341408
$ws.Cells.Item(16, 9).Value = 1.021830898894408
$ws.Cells.Item(16, 10).Value = 0.9779807578827066
$ws.Cells.Item(16, 11).Value = 1.016050029596795
$ws.Cells.Item(16, 12).Value = 0.9422283095875491
$ws.Cells.Item(16, 13).Value = 1.052557562715659
$ws.Cells.Item(16, 14).Value = 0.9648903771982843
$ws.Cells.Item(16, 15).Value = 0.9896185947905786
$ws.Cells.Item(16, 16).Value = 0.9895705468904716
